# daily auto push: 2026-02-22 03:15 UTC
# Insert one new data row at row 848 (A:D), pushing the existing rows
# 848-889 down to 849-890, and fill the new row with the 2026/02/22
# 03:00-08:00 slot (time value 8) that was missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 848 downward by inserting a new blank row above the old row 848.
$ws.Rows.Item(848).Insert()

# Populate the newly inserted row 848 with the new record. Force column A to
# remain plain text (not auto-converted to a date serial) by setting the
# number format to Text before assigning the value, matching the other
# date cells in the column which are stored as literal strings.
$ws.Cells.Item(848, 1).NumberFormat = "@"
$ws.Cells.Item(848, 1).Value = "2026/02/22"
$ws.Cells.Item(848, 2).Value = "日"
$ws.Cells.Item(848, 3).Value = 8
$ws.Cells.Item(848, 4).Value = 201
